# bug fixed for TranscriptField
# Adds a small "GPU runtime vs concurrent users" reference table (rows 13-16)
# below the existing cost table, fixes the stray border styling on E6, widens
# a couple of columns to fit the new content, and moves the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- fix E6's formatting: it had a stray leftover border style; make it match
# --- the plain styling already used by E4/E5 (copy format only, keep value)
$ws.Range("E4").Copy()
$ws.Range("E6").PasteSpecial(-4122)   # xlPasteFormats

# --- widen columns B and D so the new "concurrent"/"GPU runtime" columns read
# --- cleanly
$ws.Columns.Item(2).ColumnWidth = 9.998697916666666
$ws.Columns.Item(4).ColumnWidth = 12.166666666666666

# --- new reference table: video length / concurrent streams / GPU runtime
# header row (write C13 first so the shared-string table order matches)
$ws.Range("C13").Value = "GPU runtime"
$ws.Range("A13").Value = "video length"
$ws.Range("B13").Value = "concurrent"

# row 14: 1 concurrent stream
$ws.Range("A14").Value = 0.7368055555555556
$ws.Range("A14").NumberFormat = "h:mm"
$ws.Range("B14").Value = 1
$ws.Range("C14").Value = 0.052083333333333336
$ws.Range("C14").NumberFormat = "h:mm"

# row 15: 2 concurrent streams
$ws.Range("A15").Value = 0.7368055555555556
$ws.Range("A15").NumberFormat = "h:mm"
$ws.Range("B15").Value = 2
$ws.Range("C15").Value = 0.09930555555555555
$ws.Range("C15").NumberFormat = "h:mm"

# row 16: 3 concurrent streams
$ws.Range("A16").Value = 0.7368055555555556
$ws.Range("A16").NumberFormat = "h:mm"
$ws.Range("B16").Value = 3
$ws.Range("C16").Value = 0.13541666666666666
$ws.Range("C16").NumberFormat = "h:mm"

# --- move the active selection like the saved workbook shows
$ws.Range("H20").Select()
